# Updates the CheckIn.docx reservation data (rut, nombre, telefono,
# fechas, regalo entregado, and valor total) per "se modifico la lista
# de reserva".

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "16596395-4" "14070145-9"
Replace-Text "Pancho Melo" "Patricia Riquelme"
Replace-Text "987654321" "950700928"
Replace-Text "01-12-20" "15-11-20"
Replace-Text "08-12-20" "29-11-20"
Replace-Text "$160000" "$300000"

# entregaRegalo is the only run in the document containing exactly "No".
$rng = $d.Bookmarks("entregaRegalo").Range
$rng.Text = "Si"
